$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(2, 8, 5.8),
    @(2, 9, 6.2),
    @(2, 15, 1.29),
    @(2, 17, 1.86),
    @(2, 22, 1.19),
    @(2, 27, 170),
    @(2, 39, 120),
    @(3, 14, 1.34),
    @(3, 16, 1.34),
    @(5, 6, 5.8),
    @(5, 7, 7),
    @(5, 8, 1.54),
    @(5, 9, 1.66),
    @(5, 10, 4.3),
    @(5, 11, 5),
    @(5, 12, 1.32),
    @(5, 13, 1.05),
    @(5, 14, 4.3),
    @(5, 15, 1.24),
    @(5, 16, 2.14),
    @(5, 17, 1.7),
    @(5, 18, 1.45),
    @(5, 19, 2.74),
    @(5, 20, 1.81),
    @(5, 21, 2),
    @(5, 22, 2.52),
    @(5, 23, 1.16),
    @(5, 24, 19.5),
    @(5, 25, 12),
    @(5, 26, 10.5),
    @(5, 27, 970),
    @(5, 28, 30),
    @(5, 29, 10.5),
    @(5, 31, 16.5),
    @(5, 32, 60),
    @(5, 33, 26),
    @(5, 34, 23),
    @(5, 35, 34),
    @(5, 36, 250),
    @(5, 37, 90),
    @(5, 38, 100),
    @(5, 39, 150),
    @(5, 40, 110),
    @(5, 41, 9.800000000000001),
    @(6, 17, 1.63),
    @(6, 19, 2.4),
    @(7, 7, 2.02),
    @(7, 10, 3.3),
    @(7, 12, 1.32),
    @(7, 19, 3.1),
    @(7, 23, 1.98),
    @(8, 12, 1.36),
    @(8, 19, 3.3),
    @(8, 25, 28),
    @(8, 35, 160),
    @(8, 39, 210),
    @(8, 40, 6.8),
    @(9, 6, 2.28),
    @(9, 7, 2.54),
    @(9, 8, 3.2),
    @(9, 9, 3.8),
    @(9, 12, 1.44),
    @(9, 13, 1.08),
    @(9, 14, 3.2),
    @(9, 15, 1.36),
    @(9, 18, 1.28),
    @(9, 19, 3.75),
    @(9, 20, 1.82),
    @(9, 21, 2),
    @(9, 24, 15),
    @(9, 25, 14.5),
    @(9, 26, 29),
    @(9, 27, 80),
    @(9, 28, 11.5),
    @(9, 29, 9.199999999999999),
    @(9, 30, 18),
    @(9, 31, 55),
    @(9, 32, 18),
    @(9, 33, 14),
    @(9, 34, 23),
    @(9, 35, 70),
    @(9, 36, 36),
    @(9, 37, 34),
    @(9, 38, 55),
    @(9, 39, 140),
    @(9, 40, 29),
    @(9, 41, 60),
    @(10, 6, 3.3),
    @(10, 7, 4.9),
    @(10, 8, 1.93),
    @(10, 9, 2.46),
    @(10, 11, 5.6),
    @(10, 16, 1.66),
    @(10, 17, 1.95),
    @(10, 19, 3.45),
    @(10, 22, 1.69),
    @(11, 7, 7.4),
    @(11, 17, 2.12),
    @(11, 18, 1.31),
    @(11, 23, 1.16),
    @(11, 31, 18.5),
    @(11, 34, 27),
    @(11, 35, 46),
    @(11, 38, 130),
    @(12, 6, 2.82),
    @(12, 7, 2.84),
    @(12, 10, 3.3),
    @(12, 11, 3.35),
    @(12, 12, 1.5),
    @(12, 14, 3.1),
    @(12, 15, 1.45),
    @(12, 16, 1.71),
    @(12, 17, 2.38),
    @(12, 18, 1.26),
    @(12, 19, 4.5),
    @(12, 20, 1.97),
    @(12, 21, 1.97),
    @(12, 23, 1.54),
    @(12, 25, 9.4),
    @(12, 28, 9.4),
    @(12, 30, 13),
    @(12, 34, 20),
    @(12, 36, 44),
    @(12, 40, 36),
    @(12, 41, 40),
    @(13, 6, 4.5),
    @(13, 7, 4.8),
    @(13, 9, 2.12),
    @(13, 10, 3.3),
    @(13, 11, 3.45),
    @(13, 12, 1.49),
    @(13, 13, 1.1),
    @(13, 14, 2.98),
    @(13, 15, 1.44),
    @(13, 18, 1.25),
    @(13, 19, 4.4),
    @(13, 20, 2),
    @(13, 21, 1.87),
    @(13, 22, 1.92),
    @(13, 23, 1.27),
    @(13, 24, 11),
    @(13, 25, 7.6),
    @(13, 26, 12),
    @(13, 27, 26),
    @(13, 28, 13.5),
    @(13, 29, 7.8),
    @(13, 30, 11),
    @(13, 31, 26),
    @(13, 32, 34),
    @(13, 33, 18.5),
    @(13, 34, 22),
    @(13, 35, 48),
    @(13, 36, 140),
    @(13, 37, 70),
    @(13, 38, 100),
    @(13, 39, 180),
    @(13, 40, 120),
    @(13, 41, 21),
    @(14, 7, 3.7),
    @(14, 9, 2.54),
    @(14, 10, 3.1),
    @(14, 12, 1.57),
    @(14, 13, 1.12),
    @(14, 14, 2.68),
    @(14, 15, 1.51),
    @(14, 16, 1.57),
    @(14, 18, 1.2),
    @(14, 19, 5.1),
    @(14, 20, 2.02),
    @(14, 21, 1.84),
    @(14, 22, 1.66),
    @(14, 23, 1.37),
    @(14, 24, 9.199999999999999),
    @(14, 25, 8.4),
    @(14, 26, 15.5),
    @(14, 27, 40),
    @(14, 28, 11),
    @(14, 29, 7.2),
    @(14, 30, 13),
    @(14, 31, 36),
    @(14, 32, 26),
    @(14, 33, 16.5),
    @(14, 34, 23),
    @(14, 35, 65),
    @(14, 36, 85),
    @(14, 37, 60),
    @(14, 38, 80),
    @(14, 39, 200),
    @(14, 40, 80),
    @(14, 41, 38),
    @(15, 7, 2.16),
    @(15, 12, 1.01),
    @(15, 13, 1.01),
    @(15, 14, 2.04),
    @(15, 15, 1.26),
    @(15, 18, 1.33),
    @(15, 19, 2.6),
    @(15, 20, 1.01),
    @(15, 21, 1.01),
    @(15, 22, 1.26),
    @(15, 23, 1.86),
    @(15, 24, 25),
    @(15, 25, 23),
    @(15, 26, 44),
    @(15, 27, 100),
    @(15, 28, 15.5),
    @(15, 29, 12.5),
    @(15, 30, 24),
    @(15, 31, 65),
    @(15, 32, 20),
    @(15, 33, 15.5),
    @(15, 34, 25),
    @(15, 35, 70),
    @(15, 36, 36),
    @(15, 37, 30),
    @(15, 38, 48),
    @(15, 39, 1000),
    @(15, 40, 1000),
    @(15, 41, 1000),
    @(16, 7, 2.36),
    @(16, 8, 3.5),
    @(16, 9, 3.55),
    @(16, 12, 1.41),
    @(16, 17, 2.06),
    @(16, 20, 1.82),
    @(16, 22, 1.39),
    @(16, 23, 1.73),
    @(16, 24, 12.5),
    @(16, 26, 23),
    @(16, 27, 70),
    @(16, 28, 9.800000000000001),
    @(16, 29, 7.6),
    @(16, 30, 14),
    @(16, 31, 42),
    @(16, 32, 14),
    @(16, 34, 17.5),
    @(16, 35, 55),
    @(16, 36, 30),
    @(16, 37, 25),
    @(16, 38, 40),
    @(16, 39, 90),
    @(16, 40, 19.5),
    @(16, 41, 40),
    @(17, 6, 1.78),
    @(17, 7, 1.8),
    @(17, 8, 5.4),
    @(17, 9, 5.8),
    @(17, 11, 3.85),
    @(17, 15, 1.36),
    @(17, 24, 13.5),
    @(17, 25, 17.5),
    @(17, 27, 180),
    @(17, 35, 110),
    @(17, 36, 19),
    @(17, 37, 20),
    @(17, 39, 160),
    @(17, 40, 13),
    @(17, 41, 130),
    @(18, 7, 4.1),
    @(18, 8, 2.14),
    @(18, 10, 3.35),
    @(18, 17, 2.1),
    @(19, 6, 3.9),
    @(19, 8, 2.1),
    @(19, 14, 4.5),
    @(19, 24, 17),
    @(19, 27, 25),
    @(19, 28, 17),
    @(19, 29, 8.199999999999999),
    @(19, 33, 16.5),
    @(19, 36, 80),
    @(19, 37, 46),
    @(19, 38, 55),
    @(19, 40, 42),
    @(20, 9, 26),
    @(20, 16, 3.1)
)

foreach ($u in $updates) {
    $r = $u[0]
    $c = $u[1]
    $v = $u[2]
    $ws.Cells.Item($r, $c).Value = $v
}
